# Slide 4: flip the direction of the two arrow connectors feeding the
# "Round k" textbox (Straight Arrow Connector 88 / Straight Arrow
# Connector 90) so the arrowhead moves from the tail end to the head end.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$arrowNames = @("Straight Arrow Connector 88", "Straight Arrow Connector 90")
foreach ($name in $arrowNames) {
    $sh = $s.Shapes.Item($name)
    $sh.Line.BeginArrowheadStyle = 2   # msoArrowheadTriangle
    $sh.Line.EndArrowheadStyle = 1     # msoArrowheadNone
}
